$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 14, pushing the current row 14 (and 15) down by one.
$ws.Rows.Item(14).Insert()

# The old row 14 data is now in row 15, and the old row 15 data is now in row 16.
# Populate the brand-new row 14 with the latest week's record: same prices as the
# (now shifted) row 15, but a new date and a new volume.
$ws.Cells.Item(14, 1).Value = 11
$ws.Cells.Item(14, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(14, 3).Value = "Bíobío"
$ws.Cells.Item(14, 4).Value = 44875
$ws.Cells.Item(14, 4).Style = $ws.Cells.Item(15, 4).Style
$ws.Cells.Item(14, 4).NumberFormat = $ws.Cells.Item(15, 4).NumberFormat
$ws.Cells.Item(14, 5).Value = 8
$ws.Cells.Item(14, 6).Value = 300000000
$ws.Cells.Item(14, 7).Value = "Espárragos"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 300
$ws.Cells.Item(14, 11).Value = 1500
$ws.Cells.Item(14, 12).Value = 1600
$ws.Cells.Item(14, 13).Value = 1550
$ws.Cells.Item(14, 14).Value = "`$/kilo"
$ws.Cells.Item(14, 15).Value = "Provincia de Linares"
$ws.Cells.Item(14, 16).Value = 1550
$ws.Cells.Item(14, 17).Value = 1
$ws.Cells.Item(14, 18).Value = "Hortaliza"
